# Update "想去人数" (want-to-go) counts in column F across the four sheets.
# Values taken from the commit's regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 408
$ws.Range("F5").Value = 8481
$ws.Range("F7").Value = 10550
$ws.Range("F15").Value = 291
$ws.Range("F18").Value = 74
$ws.Range("F20").Value = 409
$ws.Range("F22").Value = 1799
$ws.Range("F24").Value = 533
$ws.Range("F25").Value = 338
$ws.Range("F26").Value = 282
$ws.Range("F30").Value = 1158
$ws.Range("F33").Value = 1413
$ws.Range("F34").Value = 436
$ws.Range("F35").Value = 341
$ws.Range("F36").Value = 283
$ws.Range("F38").Value = 127
$ws.Range("F39").Value = 510
$ws.Range("F40").Value = 341
$ws.Range("F41").Value = 91
$ws.Range("F42").Value = 281
$ws.Range("F43").Value = 633
$ws.Range("F44").Value = 29
$ws.Range("F45").Value = 88
$ws.Range("F46").Value = 87

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 46
$ws.Range("F17").Value = 379

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2795

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F9").Value = 408
$ws.Range("F10").Value = 8481
$ws.Range("F12").Value = 10550
$ws.Range("F16").Value = 291
$ws.Range("F17").Value = 74
$ws.Range("F19").Value = 1799
$ws.Range("F21").Value = 533
$ws.Range("F22").Value = 282
$ws.Range("F28").Value = 1158
$ws.Range("F34").Value = 1413
$ws.Range("F35").Value = 436
$ws.Range("F37").Value = 341
$ws.Range("F38").Value = 127
$ws.Range("F39").Value = 510
$ws.Range("F41").Value = 341
$ws.Range("F42").Value = 91
$ws.Range("F43").Value = 281
$ws.Range("F45").Value = 46
$ws.Range("F46").Value = 379
$ws.Range("F47").Value = 633
$ws.Range("F48").Value = 88
$ws.Range("F49").Value = 87
